$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 18 (shifts PULMICORT..footer down by one), to make room
# for the new product "PANTHENOL 2% TOPICAL CREAM 20 GM" which sorts
# alphabetically right after "PANADOL ADVANCE 500 MG 48 TABLETS" (row 17).
$ws.Rows("18:18").Insert()

# Populate the new row 18 with the new product's data.
$ws.Range("A18").Value = 15
$ws.Range("B18").Value = "PANTHENOL 2% TOPICAL CREAM 20 GM"
$ws.Range("H18").Value = "4:0"
$ws.Range("L18").Value = 32
$ws.Range("N18").Value = 1

# Renumber the "م" (sequence) column for all rows that shifted down.
$ws.Range("A19").Value = 16
$ws.Range("A20").Value = 17
$ws.Range("A21").Value = 18
$ws.Range("A22").Value = 19
$ws.Range("A23").Value = 20
$ws.Range("A24").Value = 21
$ws.Range("A25").Value = 22

# Update the total in the summary row (was K25, now K26).
$ws.Range("K26").Value = 1412.0999999999999
